$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 649
$ws.Range("C2").Value = 999
$ws.Range("D2").Value = 73

$ws.Range("B3").Value = 368
$ws.Range("C3").Value = 1993
$ws.Range("D3").Value = 37
$ws.Range("E3").Value = 110
